$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the start/finish dates for row 5 (story 2.1.4 / H5 & I5 were empty)
$finishedDate = Get-Date -Year 2022 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("H5").Value = $finishedDate
$ws.Range("I5").Value = $finishedDate

# Update the view: move the selection to I5 (also clears the stale topLeftCell scroll position)
$ws.Range("I5").Select()
